$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Palmer", "Butter - Salted", "1", "0.00", "0.00"),
    @("PERF", "Vegan Egg", "1", "99.59", "99.59"),
    @("", "Flour - Millers Choice", "1", "0.00", "0.00"),
    @("", "Tamper Evident - 12oz Bowl (Smoothie)", "1", "0.00", "0.00"),
    @("Web", "Bag Paper - 4.5x11.75 Window (RSS)", "1", "80.49", "80.49"),
    @("", "Sani-T-10 Sanitizer", "1", "107.36", "107.36")
)

$startRow = 32
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Force text storage (NumberFormat "@") on the numeric-looking columns
    # (Quantity / Cost Per / Total Cost) so values like "1" / "0.00" /
    # "99.59" are stored as strings, matching the rest of the sheet.
    $rng = $ws.Range("C$row" + ":E$row")
    $rng.NumberFormat = "@"

    # Column A is left blank (no value written) for rows where the source
    # SKU cell is empty, matching the original sheet's empty cells.
    if ($rowData[0] -ne "") {
        $ws.Cells.Item($row, 1).Value = $rowData[0]
    }
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
